$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (crypto price/volume updates).
# Values are written with a leading apostrophe so Excel keeps them as literal
# text (matching the workbook author's original inline-string cells) instead
# of auto-converting numeric- or percent-looking text into a number.
$updates = [ordered]@{
    "D2" = "291.37"
    "E2" = "-3.04%"
    "D3" = "30.61"
    "E3" = "-6.37%"
    "D4" = "4.945"
    "E4" = "-0.14%"
    "D5" = "0.07208"
    "E5" = "-6.30%"
    "D6" = "1.810"
    "E6" = "-7.98%"
    "D7" = "7.689"
    "E7" = "-1.81%"
    "D8" = "3.761"
    "E8" = "-1.04%"
    "D9" = "0.8963"
    "E9" = "-2.58%"
    "D10" = "0.1654"
    "E10" = "-5.61%"
    "D11" = "0.07728"
    "E11" = "-0.68%"
    "D12" = "0.08038"
    "E12" = "-6.70%"
    "D13" = "0.03036"
    "E13" = "-4.55%"
    "D14" = "0.1000"
    "E14" = "0.01%"
    "D15" = "0.001508"
    "E15" = "-0.81%"
    "D16" = "0.005738"
    "E16" = "-0.80%"
    "D18" = "3.469"
    "E18" = "0.26%"
    "E19" = "-3.31%"
    "E20" = "-0.94%"
    "D22" = "4.039"
    "E22" = "-5.59%"
    "D23" = "0.2387"
    "E23" = "19.72%"
    "D24" = "0.04510"
    "E24" = "0.00%"
    "D25" = "0.001215"
    "E25" = "-0.63%"
    "D26" = "0.004007"
    "E26" = "-9.12%"
    "D27" = "0.0001250"
    "E27" = "-0.18%"
    "D39" = "0.01582"
    "E39" = "-6.66%"
    "D40" = "0.04404"
    "E40" = "-6.03%"
    "D41" = "0.007231"
    "E41" = "-3.78%"
    "D42" = "0.009940"
    "D43" = "0.1306"
    "E43" = "-3.25%"
    "D44" = "0.002006"
    "E44" = "-14.06%"
    "D45" = "0.009502"
    "E45" = "-9.75%"
    "D46" = "0.00005990"
    "E46" = "-4.24%"
    "D47" = "0.00000000750"
    "E47" = "-0.18%"
    "E48" = "172.74%"
    "D49" = "0.002999"
    "E49" = "-3.40%"
    "D50" = "0.00002100"
    "E50" = "-0.18%"
    "D51" = "0.0002000"
    "E51" = "-0.18%"
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = "'" + $updates[$ref]
}
